$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M2").Value = "13:28"
$ws.Range("M10").Value = "14:32"
$ws.Range("M13").Value = "12:00"
$ws.Range("M14").Value = "12:00"
